$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Tag"
$ws.Range("L1").Value = "Instrument"

$ws.Range("K2").Select()
